$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 39
$ws.Range("C2").Value = 69
$ws.Range("E2").Value = 36.11111111111111
$ws.Range("F2").Value = 0.33006
$ws.Range("G2").Value = 0.023036
$ws.Range("H2").Value = 0.003688712151053987
$ws.Range("I2").Value = 0.007229875816065815
$ws.Range("J2").Value = 0.3372898758160658
$ws.Range("K2").Value = 0.3228301241839342
